$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1877.1666
$ws.Range("J17").Value = 1877.1666
$ws.Range("L17").Value = 5631.4998
$ws.Range("N17").Value = -5967.4998
$ws.Range("H18").Value = 1680.625
$ws.Range("I18").Value = 1691
$ws.Range("J18").Value = 1649.5
$ws.Range("K18").Value = 1691
$ws.Range("L18").Value = 1649.5
$ws.Range("M18").Value = -1407
$ws.Range("N18").Value = -2217.5
$ws.Range("H32").Value = 2965.9092
$ws.Range("J32").Value = 3979.5
$ws.Range("L32").Value = 3979.5
$ws.Range("N32").Value = -4631.5
$ws.Range("H49").Value = 1000
$ws.Range("I49").Value = 1000
$ws.Range("K49").Value = 3000
$ws.Range("M49").Value = -2864
$ws.Range("H55").Value = 591.8333
$ws.Range("I55").Value = 410.2
$ws.Range("K55").Value = 410.2
$ws.Range("M55").Value = -196.2
$ws.Range("H64").Value = 4474.5
$ws.Range("J64").Value = 4450
$ws.Range("L64").Value = 4450
$ws.Range("N64").Value = -4946
$ws.Range("H67").Value = 4474.5
$ws.Range("J67").Value = 4450
$ws.Range("L67").Value = 4450
$ws.Range("N67").Value = -6166
$ws.Range("H74").Value = 142041.73
$ws.Range("I74").Value = 153495.9
$ws.Range("K74").Value = 153495.9
$ws.Range("M74").Value = -152559.9
$ws.Range("H77").Value = 142041.73
$ws.Range("I77").Value = 153495.9
$ws.Range("K77").Value = 767479.5
$ws.Range("M77").Value = -762799.5
$ws.Range("H99").Value = 1304
$ws.Range("I99").Value = 289.5
$ws.Range("K99").Value = 868.5
$ws.Range("M99").Value = 629.5
$ws.Range("H113").Value = 2224.0908
$ws.Range("I113").Value = 1996.7142
$ws.Range("K113").Value = 1996.7142
$ws.Range("M113").Value = 1257.2858
$ws.Range("H137").Value = 500
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 500
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 1500
$ws.Range("M137").ClearContents()
$ws.Range("N137").Value = -6600

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H19").Value = 2725
$ws.Range("I19").Value = 600
$ws.Range("K19").Value = 600
$ws.Range("M19").Value = -371
$ws.Range("H63").Value = 5038.7144
$ws.Range("I63").Value = 4934
$ws.Range("K63").Value = 4934
$ws.Range("M63").Value = -4248
$ws.Range("H66").Value = 5038.7144
$ws.Range("I66").Value = 4934
$ws.Range("K66").Value = 24670
$ws.Range("M66").Value = -21238
$ws.Range("H88").Value = 0
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("M88").ClearContents()
$ws.Range("N88").ClearContents()
$ws.Range("H91").Value = 0
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("M91").ClearContents()
$ws.Range("N91").ClearContents()
$ws.Range("H111").Value = 29990
$ws.Range("J111").Value = 29990
$ws.Range("L111").Value = 29990
$ws.Range("N111").Value = -38170
$ws.Range("H132").Value = 2383.55
$ws.Range("I132").Value = 2287.9443
$ws.Range("K132").Value = 6863.8329
$ws.Range("M132").Value = -4333.8329

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3250
$ws.Range("I86").Value = 3250
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 3250
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -2127
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 3250
$ws.Range("I89").Value = 3250
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 16250
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -10634
$ws.Range("N89").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1513.4736
$ws.Range("I31").Value = 1383.8
$ws.Range("K31").Value = 1383.8
$ws.Range("M31").Value = -1088.8
$ws.Range("H34").Value = 1513.4736
$ws.Range("I34").Value = 1383.8
$ws.Range("K34").Value = 1383.8
$ws.Range("M34").Value = -1181.8
$ws.Range("H86").Value = 9749
$ws.Range("I86").Value = 8999
$ws.Range("J86").Value = 9999
$ws.Range("K86").Value = 8999
$ws.Range("L86").Value = 9999
$ws.Range("M86").Value = -7876
$ws.Range("N86").Value = -12245
$ws.Range("H89").Value = 9749
$ws.Range("I89").Value = 8999
$ws.Range("J89").Value = 9999
$ws.Range("K89").Value = 44995
$ws.Range("L89").Value = 49995
$ws.Range("M89").Value = -39379
$ws.Range("N89").Value = -61227
$ws.Range("H99").Value = 2123.5
$ws.Range("I99").Value = 1164.6666
$ws.Range("K99").Value = 1164.6666
$ws.Range("M99").Value = 333.3334
$ws.Range("H105").Value = 3663
$ws.Range("J105").Value = 4328.5713
$ws.Range("L105").Value = 4328.5713
$ws.Range("N105").Value = -7822.5713
$ws.Range("H126").Value = 2123.5
$ws.Range("I126").Value = 1164.6666
$ws.Range("K126").Value = 3493.9998
$ws.Range("M126").Value = -1023.9998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 9343
$ws.Range("I3").Value = 9343
$ws.Range("K3").Value = 28029
$ws.Range("M3").Value = -27917
$ws.Range("H12").Value = 574.8570999999999
$ws.Range("I12").Value = 581.75
$ws.Range("K12").Value = 1745.25
$ws.Range("M12").Value = -1572.25
$ws.Range("H64").Value = 7225
$ws.Range("I64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("M64").ClearContents()
$ws.Range("H67").Value = 7225
$ws.Range("I67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("M67").ClearContents()
$ws.Range("H109").Value = 2475
$ws.Range("I109").Value = 2475
$ws.Range("J109").Value = 0
$ws.Range("K109").Value = 7425
$ws.Range("L109").Value = 0
$ws.Range("M109").Value = -6385
$ws.Range("N109").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 16500
$ws.Range("I10").Value = 15000
$ws.Range("J10").Value = 18000
$ws.Range("K10").Value = 15000
$ws.Range("L10").Value = 18000
$ws.Range("M10").Value = -14831
$ws.Range("N10").Value = -18338
$ws.Range("H102").Value = 3249.1765
$ws.Range("I102").Value = 3228.1333
$ws.Range("K102").Value = 3228.1333
$ws.Range("M102").Value = -1606.1333
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("H122").Value = 1909.2354
$ws.Range("I122").Value = 1546.8462
$ws.Range("K122").Value = 4640.5386
$ws.Range("M122").Value = -2190.5386

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6107.115
$ws.Range("I7").Value = 2499.375
$ws.Range("K7").Value = 2499.375
$ws.Range("M7").Value = -2387.375
$ws.Range("H16").Value = 819.8
$ws.Range("I16").Value = 819.8
$ws.Range("K16").Value = 819.8
$ws.Range("M16").Value = -649.8
$ws.Range("H40").Value = 2127.4666
$ws.Range("J40").Value = 2050
$ws.Range("L40").Value = 2050
$ws.Range("N40").Value = -2322
$ws.Range("H46").Value = 3909.923
$ws.Range("J46").Value = 4594.3335
$ws.Range("L46").Value = 4594.3335
$ws.Range("N46").Value = -4970.3335
$ws.Range("H82").Value = 1052.7142
$ws.Range("I82").Value = 696
$ws.Range("K82").Value = 696
$ws.Range("M82").Value = -335
$ws.Range("H85").Value = 1052.7142
$ws.Range("I85").Value = 696
$ws.Range("K85").Value = 696
$ws.Range("M85").Value = 448.25
$ws.Range("H122").Value = 5355.2
$ws.Range("I122").Value = 4404.9414
$ws.Range("J122").Value = 7374.5
$ws.Range("K122").Value = 13214.8242
$ws.Range("L122").Value = 22123.5
$ws.Range("M122").Value = -10764.8242
$ws.Range("N122").Value = -27023.5
$ws.Range("H126").Value = 6107.115
$ws.Range("I126").Value = 2499.375
$ws.Range("K126").Value = 7498.125
$ws.Range("M126").Value = -5028.125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 49156.9
$ws.Range("J46").Value = 49156.9
$ws.Range("L46").Value = 49156.9
$ws.Range("N46").Value = -49618.9
$ws.Range("H126").Value = 1357.7273
$ws.Range("I126").Value = 1430.125
$ws.Range("J126").Value = 1164.6666
$ws.Range("K126").Value = 4290.375
$ws.Range("L126").Value = 3493.9998
$ws.Range("M126").Value = -1820.375
$ws.Range("N126").Value = -8433.9998
$ws.Range("H134").Value = 49156.9
$ws.Range("J134").Value = 49156.9
$ws.Range("L134").Value = 147470.7
$ws.Range("N134").Value = -152540.7
$ws.Range("H136").Value = 3654.1035
$ws.Range("J136").Value = 1794
$ws.Range("L136").Value = 5382
$ws.Range("N136").Value = -10482
